# Foil/isotope info for the "xs class" foils: fill in the previously blank
# rows 32-44 on Sheet1 with the isotopes produced in the Cu, Ni, Ti and Sc/V
# foils, matching the layout already used for the other foils further up
# the sheet (a bold header row of shared-string labels across B:K, with the
# isotope list running down column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cu foil isotopes (A32:A34)
$ws.Range("A32").Value = "60CU"
$ws.Range("A33").Value = "61CU"
$ws.Range("A34").Value = "64CU"

# Ni foil isotopes (A35:A37)
$ws.Range("A35").Value = "56NI"
$ws.Range("A36").Value = "57NI"
$ws.Range("A37").Value = "65NI"

# Fe isotope seen in one of the foils (A38)
$ws.Range("A38").Value = "59FE"

# Ti foil header row (B40:K40) - mirrors the other foil header rows
# (e.g. B24:F24 for Ni01..Ni05) further up the sheet, bold like them too.
$ws.Range("B40").Value = "Ti01"
$ws.Range("C40").Value = "Ti02"
$ws.Range("D40").Value = "Ti03"
$ws.Range("E40").Value = "Ti04"
$ws.Range("F40").Value = "Ti05"
$ws.Range("G40").Value = "Ti06"
$ws.Range("H40").Value = "Ti08"
$ws.Range("I40").Value = "Ti09"
$ws.Range("J40").Value = "Ti10"
$ws.Range("K40").Value = "Ti11"
$ws.Range("B40:K40").Font.Bold = $true

# Sc/V isotopes produced in the Ti foils (A41:A44)
$ws.Range("A41").Value = "46SC"
$ws.Range("A42").Value = "47SC"
$ws.Range("A43").Value = "48SC"
$ws.Range("A44").Value = "48V"

# Scroll/selection ends up further down the sheet after adding the new rows.
$ws.Range("A45").Select()
